# cryptos.xlsx refresh: updated Price (D) / Volume(1h) (E) figures for every
# coin row, plus a handful of rows whose rank (and therefore Coin/Link/Price/
# Volume) swapped with its neighbour. Row numbers below match the sheet.
#
# Several Price values are plain digits-and-one-dot (e.g. "214.51"), which
# Excel auto-parses as a number when assigned through .Value, truncating the
# original plain-text formatting. We reproduce what a user gets when typing
# a leading apostrophe (e.g. '214.51) to force text entry for those, leaving
# the already-unambiguous text values (e.g. "89.996.00" or the padded % change
# strings in column E) as plain assignments.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '89.996.00'
$ws.Range("E2").Value = '  +2.47%  '

# Row 3
$ws.Range("D3").Value = '3.207.69'
$ws.Range("E3").Value = '  -1.04%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = '''214.51'
$ws.Range("E5").Value = '  +3.54%  '

# Row 6
$ws.Range("D6").Value = '''621.51'
$ws.Range("E6").Value = '  +0.64%  '

# Row 7
$ws.Range("D7").Value = '''0.389'
$ws.Range("E7").Value = '  +0.99%  '

# Row 8
$ws.Range("D8").Value = '''0.699'
$ws.Range("E8").Value = '  +1.60%  '

# Row 9
$ws.Range("D9").Value = '''0.999'
$ws.Range("E9").Value = '  +0.02%  '

# Row 10
$ws.Range("D10").Value = '3.206.89'
$ws.Range("E10").Value = '  -0.94%  '

# Row 11
$ws.Range("D11").Value = '''0.575'
$ws.Range("E11").Value = '  +4.27%  '

# Row 12
$ws.Range("E12").Value = '  -3.03%  '

# Row 13
$ws.Range("D13").Value = '''0.0000255'
$ws.Range("E13").Value = '  +1.36%  '

# Row 14
$ws.Range("D14").Value = '''5.38'
$ws.Range("E14").Value = '  -0.46%  '

# Row 15
$ws.Range("D15").Value = '3.798.09'
$ws.Range("E15").Value = '  -0.93%  '

# Row 16
$ws.Range("D16").Value = '89.761.26'
$ws.Range("E16").Value = '  +2.49%  '

# Row 17
$ws.Range("D17").Value = '''32.75'
$ws.Range("E17").Value = '  -1.12%  '

# Row 18
$ws.Range("D18").Value = '3.206.15'
$ws.Range("E18").Value = '  -0.17%  '

# Row 19
$ws.Range("D19").Value = '''0.0000238'
$ws.Range("E19").Value = '  +74.61%  '

# Row 20
$ws.Range("D20").Value = '''3.37'
$ws.Range("E20").Value = '  +12.81%  '

# Row 21
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = '''13.37'
$ws.Range("E21").Value = '  -2.44%  '

# Row 22
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '''435.67'
$ws.Range("E22").Value = '  +2.20%  '

# Row 23
$ws.Range("D23").Value = '''8.56'
$ws.Range("E23").Value = '  -1.99%  '

# Row 24
$ws.Range("D24").Value = '''5.04'
$ws.Range("E24").Value = '  -2.96%  '

# Row 25
$ws.Range("B25").Value = 'Aptos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D25").Value = '''11.66'
$ws.Range("E25").Value = '  -1.68%  '

# Row 26
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = '''5.08'
$ws.Range("E26").Value = '  -3.18%  '

# Row 27
$ws.Range("D27").Value = '3.372.45'
$ws.Range("E27").Value = '  -0.86%  '

# Row 28
$ws.Range("D28").Value = '''75.46'
$ws.Range("E28").Value = '  +0.47%  '

# Row 29
$ws.Range("E29").Value = '  +0.08%  '

# Row 30
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '''0.998'
$ws.Range("E30").Value = '  -0.10%  '

# Row 31
$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").Value = '''0.153'
$ws.Range("E31").Value = '  -14.30%  '

# Row 32
$ws.Range("D32").Value = '''4.12'
$ws.Range("E32").Value = '  +35.21%  '

# Row 33
$ws.Range("D33").Value = '''8.41'
$ws.Range("E33").Value = '  -1.64%  '

# Row 34
$ws.Range("D34").Value = '''533.29'
$ws.Range("E34").Value = '  -2.96%  '

# Row 35
$ws.Range("B35").Value = 'PancakeSwap'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D35").Value = '''1.89'
$ws.Range("E35").Value = '  -1.04%  '

# Row 36
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").Value = '''6.86'
$ws.Range("E36").Value = '  +0.70%  '

# Row 37
$ws.Range("E37").Value = '  -2.29%  '

# Row 38
$ws.Range("D38").Value = '''22.32'
$ws.Range("E38").Value = '  -0.95%  '

# Row 39
$ws.Range("D39").Value = '''22.28'
$ws.Range("E39").Value = '  +1.95%  '

# Row 40
$ws.Range("D40").Value = '''0.998'
$ws.Range("E40").Value = '  -0.06%  '

# Row 41
$ws.Range("D41").Value = '''0.126'
$ws.Range("E41").Value = '  -7.49%  '

# Row 42
$ws.Range("E42").Value = '  +0.03%  '

# Row 43
$ws.Range("D43").Value = '''1.92'
$ws.Range("E43").Value = '  -0.80%  '

# Row 44
$ws.Range("D44").Value = '''0.372'
$ws.Range("E44").Value = '  -4.72%  '

# Row 45
$ws.Range("D45").Value = '''150.86'
$ws.Range("E45").Value = '  +0.48%  '

# Row 46
$ws.Range("D46").Value = '''171.53'
$ws.Range("E46").Value = '  -2.83%  '

# Row 47
$ws.Range("D47").Value = '''43.08'
$ws.Range("E47").Value = '  -2.10%  '

# Row 48
$ws.Range("D48").Value = '''0.123'
$ws.Range("E48").Value = '  -5.83%  '

# Row 49
$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").Value = '''1.24'
$ws.Range("E49").Value = '  -5.45%  '

# Row 50
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '''0.732'
$ws.Range("E50").Value = '  +2.84%  '

# Row 51
$ws.Range("D51").Value = '''0.616'
$ws.Range("E51").Value = '  +0.18%  '
